$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Enterprises density (per 1000 people)" row (row 13): SMEs 3.7 -> 3.65, MSMEs 30.5 -> 30.45
$ws.Range("C13").Value = "'3.65"
$ws.Range("D13").Value = "'30.45"

# "Employment (% of total)" row (row 14): SMEs 45.2 -> 45.21, MSMEs 70 -> 70.01
$ws.Range("C14").Value = "'45.21"
$ws.Range("D14").Value = "'70.01"

# "Enterprises (% of total)" row (row 16): Micro 87.7 -> 87.72, SMEs 12 -> 11.96, MSMEs 99.7 -> 99.68
$ws.Range("B16").Value = "'87.72"
$ws.Range("C16").Value = "'11.96"
$ws.Range("D16").Value = "'99.68"
